# Update "想去人数" (interested-count) figures in column F across all four
# sheets of the 广州-漫展信息 workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 873
$ws.Cells.Item(7, 6).Value = 457
$ws.Cells.Item(9, 6).Value = 2126
$ws.Cells.Item(10, 6).Value = 610
$ws.Cells.Item(11, 6).Value = 274
$ws.Cells.Item(12, 6).Value = 111
$ws.Cells.Item(13, 6).Value = 1029
$ws.Cells.Item(14, 6).Value = 171
$ws.Cells.Item(15, 6).Value = 2171
$ws.Cells.Item(16, 6).Value = 626
$ws.Cells.Item(17, 6).Value = 11585
$ws.Cells.Item(18, 6).Value = 1180
$ws.Cells.Item(19, 6).Value = 551
$ws.Cells.Item(20, 6).Value = 119
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(24, 6).Value = 256
$ws.Cells.Item(27, 6).Value = 8

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 15
$ws.Cells.Item(9, 6).Value = 146
$ws.Cells.Item(12, 6).Value = 55

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5673
$ws.Cells.Item(3, 6).Value = 470
$ws.Cells.Item(4, 6).Value = 454

# Sheet "全部类型" (All types - aggregate)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5673
$ws.Cells.Item(4, 6).Value = 470
$ws.Cells.Item(5, 6).Value = 454
$ws.Cells.Item(9, 6).Value = 873
$ws.Cells.Item(11, 6).Value = 457
$ws.Cells.Item(12, 6).Value = 15
$ws.Cells.Item(13, 6).Value = 2126
$ws.Cells.Item(14, 6).Value = 610
$ws.Cells.Item(15, 6).Value = 274
$ws.Cells.Item(17, 6).Value = 111
$ws.Cells.Item(19, 6).Value = 1029
$ws.Cells.Item(21, 6).Value = 171
$ws.Cells.Item(22, 6).Value = 146
$ws.Cells.Item(24, 6).Value = 2171
$ws.Cells.Item(25, 6).Value = 626
$ws.Cells.Item(26, 6).Value = 11585
$ws.Cells.Item(28, 6).Value = 55
$ws.Cells.Item(29, 6).Value = 1180
$ws.Cells.Item(30, 6).Value = 551
$ws.Cells.Item(31, 6).Value = 119
$ws.Cells.Item(32, 6).Value = 4
$ws.Cells.Item(38, 6).Value = 256
$ws.Cells.Item(49, 6).Value = 8
